# The 'year' column (A) held Excel date serials (42444.5 / 42628.5, shown
# with a date number format) for the two data rows. The final upload
# instead stores the plain year text "2016" as a normal (General-format)
# text cell, matching how the source data was re-exported.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A2:A3")
$rng.ClearFormats()
$rng.NumberFormat = "@"
$rng.Value = "2016"
$rng.ClearFormats()
